$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column E (statuses_count becomes column F, etc.)
$ws.Columns("E:E").Insert()

# New header for inserted column
$ws.Range("E1").Value = "id"
$ws.Range("E1").Font.Bold = $true
$ws.Range("E1").HorizontalAlignment = -4108
$ws.Range("E1").VerticalAlignment = -4160
$ws.Range("E1").Borders.LineStyle = 1

# Row 2 (count) - new column E value
$ws.Range("E2").Value = 1200

# Row 3 (mean)
$ws.Range("B3").Value = 10545.64333333333
$ws.Range("C3").Value = 337314.1641666667
$ws.Range("D3").Value = 4271.750833333334
$ws.Range("E3").Value = 117659042009188992
$ws.Range("F3").Value = 24809.70416666667
$ws.Range("G3").Value = 3062.076666666667
$ws.Range("H3").Value = 7.888849631600142
$ws.Range("I3").Value = 3.631691511453875
$ws.Range("J3").Value = 2.709531735934661
$ws.Range("K3").Value = 9.783333333333333
$ws.Range("L3").Value = 64.5325
$ws.Range("M3").Value = 0.6331464341602853
$ws.Range("N3").Value = 9.484999999999999
$ws.Range("O3").Value = 55.565
$ws.Range("P3").Value = 5.78784525451346

# Row 4 (std)
$ws.Range("B4").Value = 25797.93519269118
$ws.Range("C4").Value = 1884400.071274768
$ws.Range("D4").Value = 51751.81873956503
$ws.Range("E4").Value = 292421216130286976
$ws.Range("F4").Value = 120484.1439924343
$ws.Range("G4").Value = 1001.046727804583
$ws.Range("H4").Value = 42.77300667365063
$ws.Range("I4").Value = 9.919048219479846
$ws.Range("J4").Value = 11.47763952068743
$ws.Range("K4").Value = 8.654891382502997
$ws.Range("L4").Value = 54.51399291254852
$ws.Range("M4").Value = 0.3534234229180773
$ws.Range("N4").Value = 7.960061378382896
$ws.Range("O4").Value = 46.97342717662045
$ws.Range("P4").Value = 7.045101187305121

# Row 5 (min)
$ws.Range("B5").Value = 0
$ws.Range("C5").Value = 0
$ws.Range("D5").Value = 0
$ws.Range("E5").Value = 418
$ws.Range("F5").Value = 0
$ws.Range("G5").Value = 488
$ws.Range("H5").Value = 0
$ws.Range("I5").Value = 0
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 1
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 0
$ws.Range("N5").Value = 1
$ws.Range("O5").Value = 1
$ws.Range("P5").Value = 1

# Row 6 (25%)
$ws.Range("B6").Value = 331.5
$ws.Range("C6").Value = 33
$ws.Range("D6").Value = 29
$ws.Range("E6").Value = 94439953.25
$ws.Range("F6").Value = 1363
$ws.Range("G6").Value = 2379.75
$ws.Range("H6").Value = 0.5013136288998358
$ws.Range("I6").Value = 0.1124325701065373
$ws.Range("J6").Value = 0.0007312739345898219
$ws.Range("K6").Value = 1.75
$ws.Range("L6").Value = 12
$ws.Range("M6").Value = 0.3137138574907887
$ws.Range("N6").Value = 2
$ws.Range("O6").Value = 10.75
$ws.Range("P6").Value = 3.657894736842105

# Row 7 (50%)
$ws.Range("B7").Value = 2042.5
$ws.Range("C7").Value = 350
$ws.Range("D7").Value = 271.5
$ws.Range("E7").Value = 377230939
$ws.Range("F7").Value = 4436.5
$ws.Range("G7").Value = 3254.5
$ws.Range("H7").Value = 1.443625841750842
$ws.Range("I7").Value = 0.6632034452872573
$ws.Range("J7").Value = 0.4653561037105763
$ws.Range("K7").Value = 8
$ws.Range("L7").Value = 54
$ws.Range("M7").Value = 0.6824302202821855
$ws.Range("N7").Value = 8
$ws.Range("O7").Value = 47
$ws.Range("P7").Value = 5.192307692307693

# Row 8 (75%)
$ws.Range("B8").Value = 8328.25
$ws.Range("C8").Value = 13282.25
$ws.Range("D8").Value = 813
$ws.Range("E8").Value = 2338591177.25
$ws.Range("F8").Value = 16481.5
$ws.Range("G8").Value = 3911
$ws.Range("H8").Value = 5.358665313230251
$ws.Range("I8").Value = 2.585386326599655
$ws.Range("J8").Value = 2.187619047619048
$ws.Range("K8").Value = 16
$ws.Range("L8").Value = 113
$ws.Range("M8").Value = 0.9992692604804158
$ws.Range("N8").Value = 16
$ws.Range("O8").Value = 98
$ws.Range("P8").Value = 6.631578947368421

# Row 9 (max)
$ws.Range("B9").Value = 354763
$ws.Range("C9").Value = 41478966
$ws.Range("D9").Value = 1473166
$ws.Range("E9").Value = 1118951206448180992
$ws.Range("F9").Value = 2771910
$ws.Range("G9").Value = 5158
$ws.Range("H9").Value = 1062.034482758621
$ws.Range("I9").Value = 191.3500539374326
$ws.Range("J9").Value = 242.3333333333333
$ws.Range("K9").Value = 90
$ws.Range("L9").Value = 165
$ws.Range("M9").Value = 1
$ws.Range("N9").Value = 32
$ws.Range("O9").Value = 156
$ws.Range("P9").Value = 110
